$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("National GE")

$ws.Range("B4").Value = 0.00800554476769829
$ws.Range("C4").Value = 0.03718093429243172
$ws.Range("D4").Value = 0.01651556812775156
$ws.Range("E4").Value = 0.03014806460128146
$ws.Range("F4").Value = 0.002456945691160374
$ws.Range("G4").Value = 0.0009680266521925276
$ws.Range("H4").Value = 0.001630252534372018
$ws.Range("I4").Value = 0.003338912214397779
$ws.Range("J4").Value = 0.00613332502259482
$ws.Range("K4").Value = 0.0006783218714613196
$ws.Range("L4").Value = 0.001492421298095608
$ws.Range("M4").Value = 0.004607991906929692
$ws.Range("N4").Value = 0.0006073327254271114
$ws.Range("O4").Value = 0.001036342227098777
$ws.Range("P4").Value = 0.003134823020978324
$ws.Range("Q4").Value = 0.001047416843595187
$ws.Range("R4").Value = 0.0008243650161362194
$ws.Range("S4").Value = 0.00008984214355762234
$ws.Range("T4").Value = 0.0006088989897662932
$ws.Range("U4").Value = 0.001489625437887362
$ws.Range("V4").Value = 0.005225410486363665
$ws.Range("W4").Value = 0.00439365509946645
$ws.Range("X4").Value = 0.007688805633126627
$ws.Range("Y4").Value = 0.004057752348107797
$ws.Range("Z4").Value = 0.005297996715271988
$ws.Range("AA4").Value = 0.01002773407981162
$ws.Range("AB4").Value = 0.003353085672375054
$ws.Range("AC4").Value = 0.0006008942301709348
$ws.Range("AD4").Value = 0.002139108671080252
$ws.Range("AE4").Value = 0.004786818264616974
$ws.Range("AF4").Value = 0.01533516290812533
$ws.Range("AG4").Value = 0.009190290925008008
$ws.Range("AH4").Value = 0.002261397239270025
$ws.Range("AI4").Value = 0.005001381294806605
$ws.Range("AJ4").Value = 0.003155749554641522
